# laporan_kegiatan.docx: replace the old six-section "LAPORAN PELAKSANAAN
# KEGIATAN" body with the new shorter "LAPORAN KEGIATAN" layout (title +
# satker line, a 2-column info table, three Heading2 sections and a PPK
# signature block).

$d = $word.ActiveDocument

# --- 1. Strip the whole existing body down to a single paragraph -----------
# Word never lets you delete the very last paragraph mark of the body (it
# anchors the sectPr), so delete paragraphs 1..N-1 and keep the Nth as a
# throw-away insertion point / stand-in for the document's final line.
$n = $d.Paragraphs.Count
if ($n -gt 1) {
    $deleteRange = $d.Range($d.Paragraphs.Item(1).Range.Start, $d.Paragraphs.Item($n - 1).Range.End)
    $deleteRange.Delete()
}

# --- 2. Insert the new report body as OOXML --------------------------------
# Everything except the very last "NIP: {{ppk_nip}}" line is inserted here,
# right before the one paragraph we kept.
$newBodyXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>LAPORAN KEGIATAN</w:t></w:r></w:p>
<w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:sz w:val="22"/></w:rPr><w:t>{{satker_nama}}</w:t></w:r></w:p>
<w:p/>
<w:tbl>
<w:tblPr><w:tblStyle w:val="LightGrid-Accent1"/><w:tblW w:type="auto" w:w="0"/><w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/></w:tblPr>
<w:tblGrid><w:gridCol w:w="4320"/><w:gridCol w:w="4320"/></w:tblGrid>
<w:tr>
<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Periode</w:t></w:r></w:p></w:tc>
<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>{{periode_laporan}}</w:t></w:r></w:p></w:tc>
</w:tr>
<w:tr>
<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Kegiatan</w:t></w:r></w:p></w:tc>
<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>{{nama_kegiatan}}</w:t></w:r></w:p></w:tc>
</w:tr>
<w:tr>
<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Penanggungjawab</w:t></w:r></w:p></w:tc>
<w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>{{ppk_nama}}</w:t></w:r></w:p></w:tc>
</w:tr>
</w:tbl>
<w:p/>
<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>DESKRIPSI KEGIATAN</w:t></w:r></w:p>
<w:p><w:r><w:t>{{deskripsi_kegiatan}}</w:t></w:r></w:p>
<w:p/>
<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>HASIL DAN MANFAAT</w:t></w:r></w:p>
<w:p><w:r><w:t>{{hasil_manfaat}}</w:t></w:r></w:p>
<w:p/>
<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>KESIMPULAN</w:t></w:r></w:p>
<w:p><w:r><w:t>{{kesimpulan}}</w:t></w:r></w:p>
<w:p/>
<w:p/>
<w:p><w:r><w:t>PPK</w:t></w:r></w:p>
<w:p/>
<w:p/>
<w:p><w:r><w:t>{{ppk_nama}}</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertionPoint = $d.Paragraphs.Item(1).Range
$insertionPoint.Collapse(1)
[void]$insertionPoint.InsertXML($newBodyXml)

# --- 3. Turn the one surviving original paragraph into the closing line ----
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$lastPara.Text = 'NIP: {{ppk_nip}}'

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
Write-Output ("Tables now: " + $d.Tables.Count)
